$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 and J1, reusing the existing header style (copy format from H1)
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Add the data values for row 2 (plain numbers, no special style like the other row-2 cells)
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9
